$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2021" column (M) is being added, mirroring the formatting of the
# existing "2020" column (L). Copy L2:L10 (values+formats) into M2:M10 so the
# new column picks up the same borders / fonts / number formats as column L,
# then overwrite the copied values with the new 2021 figures.
$ws.Range("L2:L10").Copy($ws.Range("M2:M10")) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("M3").Value = 2021
$ws.Range("M4").Value = 952
$ws.Range("M5").Value = 10437
$ws.Range("M6").Value = 2253
$ws.Range("M7").Value = 8184
$ws.Range("M8").Value = 14020
$ws.Range("M9").Value = 5139
$ws.Range("M10").Value = 8881

# The sheet's active selection moved to P8 as part of this edit.
$ws.Range("P8").Select() | Out-Null
